# Updates the cryptos list (coin order, prices, links, 1h-volume deltas)
# to match the latest GitHub Actions scrape.
#
# Plain `Range.Value = <numeric-looking string>` gets auto-coerced to a
# real number by Excel's COM layer (e.g. "549.31" -> 549.31 as Double),
# which would change these cells from text to numbers -- not what the
# source data (and the target workbook) represents them as. Instead, each
# cell is populated through a text-literal formula (="...") and then
# immediately frozen to a plain value via Copy / PasteSpecial(xlPasteValues).
# That keeps the result a literal string (as the original inline-string
# cells were) and leaves styles/number-formats untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $c = $ws.Range($addr)
    $c.Formula = '="' + $text.Replace('"', '""') + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue "D2" '61.144.09'
Set-TextValue "E2" '  -0.20%  '
Set-TextValue "D3" '2.376.76'
Set-TextValue "E3" '  -0.53%  '
Set-TextValue "E4" '  +0.04%  '
Set-TextValue "D5" '549.31'
Set-TextValue "E5" '  -0.09%  '
Set-TextValue "D6" '138.94'
Set-TextValue "E6" '  -2.19%  '
Set-TextValue "E7" '  -0.01%  '
Set-TextValue "D8" '0.526'
Set-TextValue "E8" '  -1.40%  '
Set-TextValue "D9" '2.377.69'
Set-TextValue "E9" '  -0.40%  '
Set-TextValue "E10" '  +2.59%  '
Set-TextValue "E11" '  +1.26%  '
Set-TextValue "D12" '5.35'
Set-TextValue "E12" '  +0.97%  '
Set-TextValue "D13" '0.349'
Set-TextValue "E13" '  +0.66%  '
Set-TextValue "D14" '25.07'
Set-TextValue "E14" '  -1.83%  '
Set-TextValue "B15" 'ShibaInu'
Set-TextValue "C15" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D15" '0.0000166'
Set-TextValue "E15" '  +0.82%  '
Set-TextValue "B16" 'WrappedBTC'
Set-TextValue "C16" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue "D16" '61.069.72'
Set-TextValue "E16" '  -0.20%  '
Set-TextValue "B17" 'WrappedEther'
Set-TextValue "C17" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D17" '2.390.45'
Set-TextValue "E17" '  +0.15%  '
Set-TextValue "B18" 'Chainlink'
Set-TextValue "C18" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D18" '10.89'
Set-TextValue "E18" '  +1.25%  '
Set-TextValue "B19" 'Polkadot'
Set-TextValue "C19" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D19" '4.16'
Set-TextValue "E19" '  +0.59%  '
Set-TextValue "B20" 'BitcoinCash'
Set-TextValue "C20" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D20" '321.79'
Set-TextValue "E20" '  +0.92%  '
Set-TextValue "B21" 'Uniswap'
Set-TextValue "C21" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D21" '6.73'
Set-TextValue "E21" '  +0.48%  '
Set-TextValue "B22" 'Dai'
Set-TextValue "C22" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D22" '1.00'
Set-TextValue "E22" '  -0.06%  '
Set-TextValue "B23" 'Litecoin'
Set-TextValue "C23" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D23" '64.42'
Set-TextValue "E23" '  +1.27%  '
Set-TextValue "B24" 'SuiNetwork'
Set-TextValue "C24" 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue "D24" '1.70'
Set-TextValue "E24" '  -11.36%  '
Set-TextValue "B25" 'Aptos'
Set-TextValue "C25" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D25" '8.43'
Set-TextValue "E25" '  +1.47%  '
Set-TextValue "B26" 'Binance-PegBSC-USD'
Set-TextValue "C26" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D26" '1.00'
Set-TextValue "E26" '  -0.03%  '
Set-TextValue "B27" 'InternetComputer(DFINITY)'
Set-TextValue "C27" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D27" '8.18'
Set-TextValue "E27" '  +0.62%  '
Set-TextValue "B28" 'Bittensor'
Set-TextValue "C28" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D28" '506.83'
Set-TextValue "E28" '  -3.97%  '
Set-TextValue "B29" 'Kaspa'
Set-TextValue "C29" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D29" '0.150'
Set-TextValue "E29" '  +3.09%  '
Set-TextValue "B30" 'PEPE'
Set-TextValue "C30" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D30" '0.0₃0888'
Set-TextValue "E30" '  -4.31%  '
Set-TextValue "B31" 'Fetch.AI'
Set-TextValue "C31" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D31" '1.38'
Set-TextValue "E31" '  -4.02%  '
Set-TextValue "B32" 'PancakeSwap'
Set-TextValue "C32" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D32" '1.83'
Set-TextValue "E32" '  -0.60%  '
Set-TextValue "B33" 'ImmutableX'
Set-TextValue "C33" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D33" '1.52'
Set-TextValue "E33" '  -3.72%  '
Set-TextValue "B34" 'FirstDigitalUSD'
Set-TextValue "C34" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D34" '0.999'
Set-TextValue "E34" '  -0.02%  '
Set-TextValue "B35" 'NEARProtocol'
Set-TextValue "C35" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D35" '4.68'
Set-TextValue "E35" '  -0.46%  '
Set-TextValue "B36" 'Stacks'
Set-TextValue "C36" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D36" '1.90'
Set-TextValue "E36" '  +3.50%  '
Set-TextValue "B37" 'RenderToken'
Set-TextValue "C37" 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue "D37" '5.42'
Set-TextValue "E37" '  -2.11%  '
Set-TextValue "B38" 'PolygonEcosystemToken'
Set-TextValue "C38" 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue "D38" '0.379'
Set-TextValue "E38" '  +0.83%  '
Set-TextValue "B39" 'EthereumClassic'
Set-TextValue "C39" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D39" '18.61'
Set-TextValue "E39" '  +2.82%  '
Set-TextValue "B40" 'Monero'
Set-TextValue "C40" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D40" '146.59'
Set-TextValue "E40" '  +5.09%  '
Set-TextValue "B41" 'USDe'
Set-TextValue "C41" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D41" '0.999'
Set-TextValue "E41" '  -0.12%  '
Set-TextValue "B42" 'OKB'
Set-TextValue "C42" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D42" '41.32'
Set-TextValue "E42" '  +2.55%  '
Set-TextValue "B43" 'Aave'
Set-TextValue "C43" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D43" '150.33'
Set-TextValue "E43" '  +6.71%  '
Set-TextValue "B44" 'Filecoin'
Set-TextValue "C44" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D44" '3.60'
Set-TextValue "E44" '  -0.86%  '
Set-TextValue "B45" 'dogwifhat'
Set-TextValue "C45" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D45" '2.12'
Set-TextValue "E45" '  -1.73%  '
Set-TextValue "B46" 'Hedera'
Set-TextValue "C46" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D46" '0.0522'
Set-TextValue "E46" '  +0.33%  '
Set-TextValue "B47" 'Mantle'
Set-TextValue "C47" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D47" '0.577'
Set-TextValue "E47" '  +0.11%  '
Set-TextValue "B48" 'InjectiveProtocol'
Set-TextValue "C48" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D48" '19.29'
Set-TextValue "E48" '  -4.17%  '
Set-TextValue "B49" 'Stellar'
Set-TextValue "C49" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D49" '0.0911'
Set-TextValue "E49" '  +0.48%  '
Set-TextValue "B50" 'VeChain'
Set-TextValue "C50" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D50" '0.0224'
Set-TextValue "E50" '  -1.00%  '
Set-TextValue "B51" 'WhiteBITCoin'
Set-TextValue "C51" 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue "D51" '11.42'
Set-TextValue "E51" '  +0.38%  '

$excel.CutCopyMode = $false
